$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet/tab to reflect the new "through" date.
$ws.Name = "Through 2021-10-01"

# 2. "September (through 09-30)" is now a complete month -> just "September".
$ws.Range("A10").Value = "September"

# 3. Push the "Total" row (row 11) down to row 12, copying both its values
#    and its formatting (bold / centered / bordered label style) so the new
#    row 12 looks exactly like the old row 11 did.
$ws.Range("A11:H11").Copy($ws.Range("A12:H12"))
$excel.CutCopyMode = $false

# 4. Row 11 now becomes the new "October (through 10-01)" row. It already
#    has the right label style from the copy above, so just replace values.
$ws.Range("A11").Value = "October (through 10-01)"
$ws.Range("B11").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("D11").Value = 6
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 4
$ws.Range("H11").Value = 8

# 5. Update the (now shifted-down) Total row with the new column sums.
$ws.Range("B12").Value = 226
$ws.Range("C12").Value = 429
$ws.Range("D12").Value = 633
$ws.Range("E12").Value = 551
$ws.Range("F12").Value = 423
$ws.Range("G12").Value = 905
$ws.Range("H12").Value = 1256

# 6. Narrow column A slightly to fit the new, shorter month labels.
$ws.Columns.Item(1).ColumnWidth = 22.8
